$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 14 ("Alternativa..." row),
# shifting it (and the two rows below it) down by one. This also pushes
# the B14:B16 merged cell down to B15:B17 automatically.
$ws.Rows("14").Insert()

# Copy the formatting (borders/fill/font/alignment) of the row above
# (the last row of the B6:B13 block) onto the freshly inserted, still
# blank row so its cells get styles matching the surrounding table.
$srcFormat = $ws.Range("B13:D13")
$dstFormat = $ws.Range("B14:D14")
$srcFormat.Copy()
$dstFormat.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Re-extend the B6:B13 merged "Cenário Normal" cell down to the new row 14.
$ws.Range("B6:B14").Merge()

# Replace the old step-6 text with the new wording, now living in D12.
$ws.Range("D12").Value = "6. Remove carro da lista de produção"

# Match the recorded UI selection from the diff.
$ws.Range("C10").Select()
